# Automatische test-sync: 2025-08-03 18:26:50
#
# Adds a new "Testmail #7" row to the Logs sheet (row 35), extends the
# conditional-formatting ranges that covered rows 2-34 to now cover
# rows 2-35, and swaps the order of the "Intern verzoek / Actie voor
# medewerker" / "Inkoop / Bestellingen" rows on the Dashboard sheet,
# updating their counts.

$wb = $excel.ActiveWorkbook

# --- Logs sheet: append the new row ------------------------------------
$logs = $wb.Worksheets.Item("Logs")

$newRow = 35
$logs.Cells.Item($newRow, 1).Value = "Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 2).Value = "mailmind.test@zohomail.eu"
$logs.Cells.Item($newRow, 3).Value = "Testmail #7: Is dit artikel nog op voorraad?"
$logs.Cells.Item($newRow, 4).Value = "Inkoop / Bestellingen"
$logs.Cells.Item($newRow, 5).Value = "Bedankt, we hebben dit doorgestuurd naar inkoop@bedrijf.nl."
$logs.Cells.Item($newRow, 6).Value = "2025-08-03 18:26:21"
$logs.Cells.Item($newRow, 7).Value = "Ja"
$logs.Cells.Item($newRow, 8).Value = "Ja"
$logs.Cells.Item($newRow, 9).Value = "Nee"
$logs.Cells.Item($newRow, 10).Value = "Nee"

# --- Logs sheet: stretch the conditional formatting ranges to row 35 ---
foreach ($colLetter in @("D", "G", "H", "I", "J")) {
    $oldRng = $logs.Range($colLetter + "2:" + $colLetter + "34")
    $newRng = $logs.Range($colLetter + "2:" + $colLetter + "35")
    $fcs = $oldRng.FormatConditions
    for ($i = 1; $i -le $fcs.Count; $i++) {
        $fcs.Item($i).ModifyAppliesToRange($newRng)
    }
}

# --- Dashboard sheet: swap the two category rows + update counts -------
$dash = $wb.Worksheets.Item("Dashboard")

$dash.Cells.Item(4, 1).Value = "Inkoop / Bestellingen"
$dash.Cells.Item(4, 2).Value = 7
$dash.Cells.Item(5, 1).Value = "Intern verzoek / Actie voor medewerker"
$dash.Cells.Item(5, 2).Value = 6
